$wb = $excel.ActiveWorkbook

# --- 1. Rename the "Requested quantity" header on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet after the existing sheets ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Bring over the header style (bold, centered, bordered) used on the other sheets
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Bring over the date-column style (numFmt) used for column A on the other sheets
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A56").PasteSpecial(-4122)

# --- 3. Headers ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# --- 4. Data rows ---
$rows = @(
  @(44983.99999999999,4,-6.357879645996588,14.06537704020664),
  @(44997.99999999999,4,-5.720249323536112,14.91515733300269),
  @(45004.99999999999,4,-5.920047303239733,14.55482311051806),
  @(45109.99999999999,6,-5.049259802325854,15.96433572227706),
  @(45116.99999999999,6,-4.085302627984843,16.43714179096801),
  @(45123.99999999999,6,-4.903894844940293,15.59654555455128),
  @(45130.99999999999,6,-3.357883741906274,16.27294841088441),
  @(45137.99999999999,6,-3.976874213567342,16.39598217124352),
  @(45144.99999999999,6,-3.559257405099649,16.89801792808694),
  @(45151.99999999999,6,-4.036530256754497,16.64005545956788),
  @(45158.99999999999,6,-3.692911233657004,16.10293143117822),
  @(45165.99999999999,7,-3.194415785419802,16.85798881373678),
  @(45172.99999999999,7,-3.624498885689487,16.7796741177579),
  @(45179.99999999999,7,-4.098316537279386,16.86110052914108),
  @(45186.99999999999,7,-3.763317489476607,17.07684897020914),
  @(45193.99999999999,7,-3.221282488303101,16.71167305199617),
  @(45200.99999999999,7,-3.539353609679557,17.32813633151301),
  @(45207.99999999999,7,-3.159357182859559,17.72027678363642),
  @(45214.99999999999,7,-2.28669060972812,17.38404132718021),
  @(45221.99999999999,7,-2.293798670683504,17.75421391290822),
  @(45228.99999999999,7,-2.686662643139896,17.95091029292302),
  @(45235.99999999999,8,-2.039875607320992,18.76918605479878),
  @(45242.99999999999,8,-2.898621463837814,17.27429363571776),
  @(45249.99999999999,8,-2.863358174360682,17.37140617087447),
  @(45256.99999999999,8,-2.576727308458959,18.71645875668842),
  @(45263.99999999999,8,-1.084883087493013,18.54026266726752),
  @(45270.99999999999,8,-1.651405981996901,18.52007159197403),
  @(45277.99999999999,8,-2.090120890804041,18.66697548335438),
  @(45305.99999999999,9,-1.979179598225042,19.14822372394305),
  @(45319.99999999999,9,-1.247660177817743,19.757595733338),
  @(45326.99999999999,9,-1.580503238458067,19.15023753235316),
  @(45333.99999999999,9,-0.9992963415726708,19.08106858459541),
  @(45340.99999999999,9,-1.01719345886534,19.27404908634615),
  @(45347.99999999999,9,-1.309435958693907,18.94090703440994),
  @(45354.99999999999,9,-0.9458117349364114,19.47700476150857),
  @(45361.99999999999,9,-0.9648584013440996,19.33864929338958),
  @(45375.99999999999,10,-0.6321605870865095,20.45183958481082),
  @(45417.99999999999,10,-0.3297609678833329,19.77650273126551),
  @(45424.99999999999,10,0.5491557088511065,21.07255583025981),
  @(45438.99999999999,10,0.003632676435414118,20.29823567740717),
  @(45452.99999999999,11,0.5835524577639108,20.55907180943974),
  @(45487.99999999999,11,1.46739764186216,20.81738649959892),
  @(45494.99999999999,11,1.170128947927731,21.54773443504615),
  @(45522.99999999999,12,1.076250096111836,21.95051653893422),
  @(45529.99999999999,12,1.288531452773089,22.07379058670045),
  @(45536.99999999999,12,1.666091203594009,22.04745766614628),
  @(45557.99999999999,12,1.4503699522222,22.18253122312818),
  @(45564.99999999999,12,2.168346953462976,22.2524599308795),
  @(45571.99999999999,12,1.674749018113665,22.38788321506934),
  @(45578.99999999999,12,2.158200682736066,22.07938721493936),
  @(45585.99999999999,13,2.4919754457492,22.32218290896741),
  @(45592.99999999999,13,2.452730950274982,22.99851638113197),
  @(45599.99999999999,13,2.089262072741606,23.26851104951523),
  @(45606.99999999999,13,2.781016295389608,22.46099540864799),
  @(45613.99999999999,13,2.906844411850551,23.09464891554568)
)

$r = 2
foreach ($row in $rows) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
